# IKD update: GaN CMOS 2026-02-07T23:34Z
# Append 4 duplicate rows for the new "Overvoltage Suppression Filter
# Development for GaN Inverter-Fed Electrical Drive with Long Cable
# Based on Impedance Measurement" record to the Master sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

$title   = "Overvoltage Suppression Filter Development for GaN Inverter-Fed Electrical Drive with Long Cable Based on Impedance Measurement"
$year    = 2026
$pub     = "MDPI AG"
$venue   = "Electronics"
$authors = "Kroičs, Kaspars; Voitkāns, Jānis"
$doi     = "10.3390/electronics15030717"
$url     = "https://doi.org/10.3390/electronics15030717"
$doctype = "Journal"
$device  = "Inverter"
$method  = "Experiment"
$enabler = "Contacts"
$conf    = "High"
$added   = "2026-02-07"

$startRow = 155
for ($i = 0; $i -lt 4; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 2).Value  = $title      # B - Title
    $ws.Cells.Item($r, 3).Value  = $year        # C - Year
    $ws.Cells.Item($r, 4).Value  = $pub         # D - Publisher
    $ws.Cells.Item($r, 5).Value  = $venue       # E - Venue
    $ws.Cells.Item($r, 6).Value  = $authors     # F - Authors
    $ws.Cells.Item($r, 8).Value  = $doi         # H - DOI
    $ws.Cells.Item($r, 9).Value  = $url         # I - URL
    $ws.Cells.Item($r, 10).Value = $doctype     # J - DocType
    $ws.Cells.Item($r, 11).Value = $device      # K - DeviceType
    $ws.Cells.Item($r, 12).Value = $method      # L - Method
    $ws.Cells.Item($r, 13).Value = $enabler     # M - EnablerCategory
    $ws.Cells.Item($r, 17).Value = $title       # Q - EvidenceSnippet
    $ws.Cells.Item($r, 18).Value = $conf        # R - TagConfidence

    # S - AddedDate: stored as plain text ("2026-02-07"), matching every
    # other AddedDate cell in the sheet. Excel auto-converts an
    # ISO-date-shaped literal into a date serial on assignment, so force
    # the cell to Text first, then restore the Normal style so no
    # residual number-format is left behind on the cell.
    $cell = $ws.Cells.Item($r, 19)
    $cell.NumberFormat = "@"
    $cell.Value = $added
    $cell.Style = "Normal"
}
